$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E sometimes contain numeric-looking text (e.g. "548.40",
# "4.74") that must stay plain text with exact formatting (trailing zeros,
# thousands-style dots, padded percent strings). Force text entry by
# pre-formatting the range as Text, assign the literal values, then strip
# the temporary formatting back off so no stray number-format is left on
# the cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '57.698.39'
$ws.Range("E2").Value = '  -3.54%  '
$ws.Range("D3").Value = '2.920.23'
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '548.40'
$ws.Range("E5").Value = '  -2.86%  '
$ws.Range("D6").Value = '129.95'
$ws.Range("E6").Value = '  +5.44%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +1.83%  '
$ws.Range("D9").Value = '2.913.72'
$ws.Range("E9").Value = '  -1.58%  '
$ws.Range("E10").Value = '  -2.06%  '
$ws.Range("D11").Value = '4.74'
$ws.Range("E11").Value = '  -2.82%  '
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").Value = '32.57'
$ws.Range("E14").Value = '  +1.54%  '
$ws.Range("E15").Value = '  +1.70%  '
$ws.Range("D16").Value = '3.405.95'
$ws.Range("E16").Value = '  -1.60%  '
$ws.Range("E17").Value = '  +6.02%  '
$ws.Range("D18").Value = '2.920.32'
$ws.Range("E18").Value = '  -1.67%  '
$ws.Range("D19").Value = '57.701.96'
$ws.Range("E19").Value = '  -3.57%  '
$ws.Range("D20").Value = '415.29'
$ws.Range("E20").Value = '  -2.05%  '
$ws.Range("E21").Value = '  +1.62%  '
$ws.Range("D22").Value = '0.689'
$ws.Range("E22").Value = '  +3.78%  '
$ws.Range("E23").Value = '  +4.49%  '
$ws.Range("E24").Value = '  +0.21%  '
$ws.Range("D25").Value = '79.37'
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  -1.87%  '
$ws.Range("E29").Value = '  +4.19%  '
$ws.Range("D30").Value = '7.33'
$ws.Range("D31").Value = '25.12'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("E32").Value = '  -1.91%  '
$ws.Range("D33").Value = '0.0964'
$ws.Range("E33").Value = '  +1.20%  '
$ws.Range("D34").Value = '5.64'
$ws.Range("E34").Value = '  +2.17%  '
$ws.Range("D35").Value = '0.928'
$ws.Range("E35").Value = '  +1.34%  '
$ws.Range("E36").Value = '  +3.80%  '
$ws.Range("D37").Value = '0.0₃0687'
$ws.Range("E37").Value = '  +7.74%  '
$ws.Range("D38").Value = '48.10'
$ws.Range("E38").Value = '  -4.02%  '
$ws.Range("D39").Value = '8.71'
$ws.Range("E39").Value = '  +3.18%  '
$ws.Range("E40").Value = '  +8.03%  '
$ws.Range("E41").Value = '  +0.50%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.0343'
$ws.Range("E42").Value = '  -2.32%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.695.00'
$ws.Range("E43").Value = '  +1.64%  '
$ws.Range("D44").Value = '371.80'
$ws.Range("E44").Value = '  +2.63%  '
$ws.Range("D46").Value = '123.60'
$ws.Range("E46").Value = '  +2.89%  '
$ws.Range("E47").Value = '  +1.98%  '
$ws.Range("E48").Value = '  +0.74%  '
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("D50").Value = '22.70'
$ws.Range("E50").Value = '  -1.35%  '
$ws.Range("E51").Value = '  -0.07%  '

$dataRange.ClearFormats()
